$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.293.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.677.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.146"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -4.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.436"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.294.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.677.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.302.59"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("E19").Value = "  -4.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "470.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.820.44"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.06%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.04%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.995"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.653.53"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.157"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0901"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "167.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.36%  "
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000278"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "27.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.45%  "
